# Corrected excel sheets for application fix issues
#
# - Summary sheet: fee amounts in row 4 (A4, B4) bumped from 50 to 100.
# - Repayment Schedule sheet: disbursement-time fee cells (I2, K2, L2)
#   bumped from 50 to 100 to match the Summary sheet.
# - Transactions sheet: the disbursement transaction's Amount/Fees
#   cells (E2, H2) bumped from 50 to 100 to match.
# - Active sheet/selection moves from NewLoanInput -> Transactions, with
#   the cell selections on Summary and Repayment Schedule left pointing
#   at the cells the user last touched while verifying the fix.

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Value = 100
$wsSummary.Range("B4").Value = 100
$wsSummary.Range("H22").Select()

$wsRepayment = $wb.Worksheets.Item("Repayment Schedule")
$wsRepayment.Range("I2").Value = 100
$wsRepayment.Range("K2").Value = 100
$wsRepayment.Range("L2").Value = 100
$wsRepayment.Range("L2").Select()

$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("E2").Value = 100
$wsTransactions.Range("H2").Value = 100
$wsTransactions.Activate()
$wsTransactions.Range("K10").Select()
